$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(3).Delete()
$ws.Range("J12").Select()
